# Weekly price-report update: a new daily price record is inserted at the
# top of this variety's data block (row 32), pushing all the existing
# records for "Poroto verde" (Vega Modelo de Temuco) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 32 - this shifts rows 32:135 down to
# 33:136 (and the sheet's used-range dimension grows from R135 to R136
# automatically).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with this week's record.
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44676
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112031
$ws.Cells.Item(32, 7).Value = "Poroto verde"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 40
$ws.Cells.Item(32, 11).Value = 25000
$ws.Cells.Item(32, 12).Value = 25000
$ws.Cells.Item(32, 13).Value = 25000
$ws.Cells.Item(32, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 1000
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
